$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting the existing data rows down
$ws.Rows.Item(1).Insert()

# Populate the new header row
$ws.Range("A1").Value = "Transaction_CUI"
$ws.Range("B1").Value = "Status"
